$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Grow the log table: duplicate the "day separator" row (94, TUESDAY-style)
# and the last data row (96) as format templates for the new rows 98-104 ---
$ws.Range("A94:F94").Copy($ws.Range("A98:F98"))
$ws.Range("A96:F96").Copy($ws.Range("A99:F104"))

# Row 99's last column (F) uses the wrap-only style (no custom number format),
# same as e.g. F7, instead of the s=25 style copied from row96 -- fix the
# format only (keep everything else), then set the value afterwards.
$ws.Range("F7").Copy()
$ws.Range("F99").PasteSpecial(-4122)

# --- Row 98: day-of-week separator ---
$ws.Range("C98").Value = "TUESDAY"

# --- Row 99 ---
$ws.Range("A99").Value = "Pickup Mic"
$ws.Range("B99").Value = 42613
$ws.Range("C99").Value = "1630"
$ws.Range("D99").Value = "FC"
$ws.Range("E99").Value = "152"
$ws.Range("F99").Value = "Leave mic cables in place. Remove mic, stand, remote and wireless keyboard. Turn off all else as usual. To Founders 156A"

# --- Row 100 ---
$ws.Range("A100").Value = "Other"
$ws.Range("B100").Value = 42613
$ws.Range("C100").Value = "1645"
$ws.Range("D100").Value = "MC"
$ws.Range("E100").Value = "157A"
$ws.Range("F100").Value = "Door code 11012* "

# --- Row 101 ---
$ws.Range("A101").Value = "Other"
$ws.Range("B101").Value = 42613
$ws.Range("C101").Value = "1645"
$ws.Range("D101").Value = "MC"
$ws.Range("E101").Value = "157B"
$ws.Range("F101").Value = "Door code 10077* "

# --- Row 102 ---
$ws.Range("A102").Value = "Pickup PC"
$ws.Range("B102").Value = 42613
$ws.Range("C102").Value = "1700"
$ws.Range("D102").Value = "VC"
$ws.Range("E102").Value = "113"
$ws.Range("F102").Value = "Flat screen DLP and wireless keyboard. To Vanier 132 storeroom."

# --- Row 103 ---
$ws.Range("A103").Value = "Setup PC"
$ws.Range("B103").Value = 42613
$ws.Range("C103").Value = "1745"
$ws.Range("D103").Value = "MC"
$ws.Range("E103").Value = "140-SCR"
$ws.Range("F103").Value = "Portable screen is there. Door code 7083* "

# --- Row 104 (slightly taller custom row height, like other "last of day" rows) ---
$ws.Range("A104").Value = "Pickup PC"
$ws.Range("B104").Value = 42613
$ws.Range("C104").Value = "2130"
$ws.Range("D104").Value = "MC"
$ws.Range("E104").Value = "140-SCR"
$ws.Range("F104").Value = "Door code 7083*  Leave portable screen. Return to Founders 156A."
$ws.Rows(104).RowHeight = 15

# --- Extend the Task_type / Building dropdown validations down to the new rows ---
$ws.Range("A104:A1048576").Validation.Delete()
$ws.Range("A104:A1048576").Validation.Add(3, 1, 1, "=Task_type")

$ws.Range("D104:D1048576").Validation.Delete()
$ws.Range("D104:D1048576").Validation.Add(3, 1, 1, "=Building")

# --- Update frozen-pane scroll position & selection to match where the
# user ended up after entering the new rows ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 80
$ws.Range("C107").Select()
